# Scheduled-runner style refresh of the per-profession market/profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit* columns H:N) on
# each leveling-sheet Table. Values are plain numeric overwrites (no
# formulas anywhere in these sheets) coming from an external price feed;
# a handful of rows lose/gain their trailing LeveProfitHQ (N) or
# LeveProfitNQ (M) cell entirely when the HQ/NQ price recipe disappears
# or (re)appears for that leve.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1217.6842
$ws.Range("I2").Value = 150.11111
$ws.Range("J2").Value = 2178.5
$ws.Range("K2").Value = 150.11111
$ws.Range("L2").Value = 2178.5
$ws.Range("M2").Value = -37.11111
$ws.Range("N2").Value = -2404.5
$ws.Range("H9").Value = 254.85715
$ws.Range("I9").Value = 297.14285
$ws.Range("J9").Value = 212.57143
$ws.Range("K9").Value = 297.14285
$ws.Range("L9").Value = 212.57143
$ws.Range("M9").Value = -128.14285
$ws.Range("N9").Value = -550.57143
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
$ws.Range("H32").Value = 2996.3333
$ws.Range("J32").Value = 2996.3333
$ws.Range("L32").Value = 2996.3333
$ws.Range("N32").Value = -3648.3333
$ws.Range("H58").Value = 1879.1428
$ws.Range("I58").Value = 525.6667
$ws.Range("K58").Value = 1577.0001
$ws.Range("M58").Value = -1427.0001
$ws.Range("H64").Value = 4865.5
$ws.Range("J64").Value = 5209.25
$ws.Range("L64").Value = 5209.25
$ws.Range("N64").Value = -5705.25
$ws.Range("H67").Value = 4865.5
$ws.Range("J67").Value = 5209.25
$ws.Range("L67").Value = 5209.25
$ws.Range("N67").Value = -6925.25
$ws.Range("H86").Value = 3997.8572
$ws.Range("I86").Value = 4166.3335
$ws.Range("J86").Value = 3871.5
$ws.Range("K86").Value = 4166.3335
$ws.Range("L86").Value = 3871.5
$ws.Range("M86").Value = -3043.3335
$ws.Range("N86").Value = -6117.5
$ws.Range("H89").Value = 3997.8572
$ws.Range("I89").Value = 4166.3335
$ws.Range("J89").Value = 3871.5
$ws.Range("K89").Value = 20831.6675
$ws.Range("L89").Value = 19357.5
$ws.Range("M89").Value = -15215.6675
$ws.Range("N89").Value = -30589.5
$ws.Range("H112").Value = 2291.7827
$ws.Range("I112").Value = 1624.5
$ws.Range("K112").Value = 4873.5
$ws.Range("M112").Value = -3765.5
$ws.Range("H141").Value = 3103.647
$ws.Range("I141").Value = 2858.7693
$ws.Range("K141").Value = 8576.3079
$ws.Range("M141").Value = -3396.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 69999.39999999999
$ws.Range("J95").Value = 69999.39999999999
$ws.Range("L95").Value = 69999.39999999999
$ws.Range("N95").Value = -75491.39999999999
$ws.Range("H98").Value = 355
$ws.Range("J98").Value = 355
$ws.Range("L98").Value = 355
$ws.Range("N98").Value = -6345
$ws.Range("H101").Value = 39333
$ws.Range("J101").Value = 39333
$ws.Range("L101").Value = 39333
$ws.Range("N101").Value = -45823
$ws.Range("H112").Value = 19993.4
$ws.Range("J112").Value = 19993.4
$ws.Range("L112").Value = 19993.4
$ws.Range("N112").Value = -22947.4
$ws.Range("H114").Value = 39624.5
$ws.Range("J114").Value = 39624.5
$ws.Range("L114").Value = 39624.5
$ws.Range("N114").Value = -48302.5
$ws.Range("H119").Value = 64500
$ws.Range("J119").Value = 64500
$ws.Range("L119").Value = 64500
$ws.Range("N119").Value = -74176

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 1575
$ws.Range("I41").Value = 1575
$ws.Range("K41").Value = 1575
$ws.Range("M41").Value = -1147
$ws.Range("H58").Value = 55149.633
$ws.Range("I58").Value = 61049.65
$ws.Range("K58").Value = 61049.65
$ws.Range("M58").Value = -60846.65
$ws.Range("H60").Value = 27250
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 33000
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 33000
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -34022
$ws.Range("H68").Value = 54333.332
$ws.Range("J68").Value = 54333.332
$ws.Range("L68").Value = 54333.332
$ws.Range("N68").Value = -55831.332
$ws.Range("H71").Value = 54333.332
$ws.Range("J71").Value = 54333.332
$ws.Range("L71").Value = 162999.996
$ws.Range("N71").Value = -170487.996
$ws.Range("H74").Value = 49999.11
$ws.Range("J74").Value = 49999.11
$ws.Range("L74").Value = 49999.11
$ws.Range("N74").Value = -51747.11
$ws.Range("H77").Value = 49999.11
$ws.Range("J77").Value = 49999.11
$ws.Range("L77").Value = 149997.33
$ws.Range("N77").Value = -158733.33
$ws.Range("H136").Value = 55149.633
$ws.Range("I136").Value = 61049.65
$ws.Range("K136").Value = 183148.95
$ws.Range("M136").Value = -180598.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 676.8
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 6000
$ws.Range("N5").Value = -6224
$ws.Range("H37").Value = 59933.332
$ws.Range("J37").Value = 59933.332
$ws.Range("L37").Value = 179799.996
$ws.Range("N37").Value = -180023.996
$ws.Range("H39").Value = 5927.727
$ws.Range("J39").Value = 7112
$ws.Range("L39").Value = 21336
$ws.Range("N39").Value = -21924
$ws.Range("H48").Value = 2975
$ws.Range("J48").Value = 2975
$ws.Range("L48").Value = 8925
$ws.Range("N48").Value = -9425
$ws.Range("H55").Value = 8267.817999999999
$ws.Range("I55").Value = 3477
$ws.Range("J55").Value = 9332.444
$ws.Range("K55").Value = 10431
$ws.Range("L55").Value = 27997.332
$ws.Range("M55").Value = -10254
$ws.Range("N55").Value = -28351.332
$ws.Range("H129").Value = 7727.316
$ws.Range("I129").Value = 11443.6
$ws.Range("K129").Value = 34330.8
$ws.Range("M129").Value = -29330.8
$ws.Range("H135").Value = 676.8
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 25255.834
$ws.Range("I41").Value = 22886.25
$ws.Range("K41").Value = 22886.25
$ws.Range("M41").Value = -22531.25
$ws.Range("H43").Value = 7097.75
$ws.Range("J43").Value = 14999
$ws.Range("L43").Value = 14999
$ws.Range("N43").Value = -15301
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = $null
$ws.Range("H57").Value = 48999
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 48999
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 48999
$ws.Range("M57").Value = $null
$ws.Range("N57").Value = -50639
$ws.Range("H80").Value = 1359.9788
$ws.Range("J80").Value = 1772.1852
$ws.Range("L80").Value = 1772.1852
$ws.Range("N80").Value = -3768.1852
$ws.Range("H83").Value = 1359.9788
$ws.Range("J83").Value = 1772.1852
$ws.Range("L83").Value = 8860.925999999999
$ws.Range("N83").Value = -18844.926
$ws.Range("H102").Value = 3345.0908
$ws.Range("I102").Value = 3137
$ws.Range("K102").Value = 3137
$ws.Range("M102").Value = -1515
$ws.Range("H111").Value = 50003
$ws.Range("J111").Value = 50003
$ws.Range("L111").Value = 50003
$ws.Range("N111").Value = -56137

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2599
$ws.Range("I7").Value = 2599
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2599
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2487
$ws.Range("N7").Value = $null
$ws.Range("H46").Value = 7833.6
$ws.Range("J46").Value = 2680.4375
$ws.Range("L46").Value = 2680.4375
$ws.Range("N46").Value = -3056.4375
$ws.Range("H82").Value = 2450.3157
$ws.Range("I82").Value = 2140
$ws.Range("K82").Value = 2140
$ws.Range("M82").Value = -1779
$ws.Range("H85").Value = 2450.3157
$ws.Range("I85").Value = 2140
$ws.Range("K85").Value = 2140
$ws.Range("M85").Value = -892
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = $null
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null
$ws.Range("H126").Value = 2599
$ws.Range("I126").Value = 2599
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7797
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5327
$ws.Range("N126").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 31979
$ws.Range("J41").Value = 31979
$ws.Range("L41").Value = 31979
$ws.Range("N41").Value = -32759
$ws.Range("H58").Value = 18999.6
$ws.Range("J58").Value = 19999
$ws.Range("L58").Value = 19999
$ws.Range("N58").Value = -20615
$ws.Range("H100").Value = 1013.46155
$ws.Range("I100").Value = 970.4545000000001
$ws.Range("K100").Value = 1940.909
$ws.Range("M100").Value = -1399.909
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
$ws.Range("H132").Value = 43042.36
$ws.Range("I132").Value = 44798.793
$ws.Range("J132").Value = 888
$ws.Range("K132").Value = 134396.379
$ws.Range("L132").Value = 2664
$ws.Range("M132").Value = -131866.379
$ws.Range("N132").Value = -7724
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
